$d = $word.ActiveDocument

# First paragraph: add a paragraph border (5-twip spacing on all sides)
# and widen the left indent from 120 -> 225 twips.
$p1 = $d.Paragraphs(1)
$b = $p1.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5
$p1.LeftIndent = 11.25

# Replace the placeholder ID text (and drop the trailing space run) with the
# updated topic id.
$d.Content.Find.Execute("**ID__AFFARS_mp_5315_3_topic_19__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP_5315_3_3_7__ID**", 2)
